$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).ColumnWidth = 14.7109375

$values = @(
    @(-0.26762998866126964, 0.26699348880082141),
    @(-0.22088144599350557, 0.21840902312648325),
    @(-0.075634033336577033, 0.075283142809432491),
    @(-0.1372811416864721, 0.13662918772163124),
    @(-0.13062918787434175, 0.12933922137597698),
    @(-0.075524307953167558, 0.075442440315649062),
    @(-0.055442440506327628, 0.055278943329396313),
    @(-0.035278943521568351, 0.03518681731188078),
    @(-0.029186817472338866, 0.029122999708006247),
    @(-0.023122999869684691, 0.023117164065332929),
    @(-0.018617164223481808, 0.018609141039824806),
    @(-0.012609141201733287, 0.012593231023850127),
    @(-0.0065932311864820292, 0.006591007065518184),
    @(-0.027085958538181742, 0.027053447031449629),
    @(-0.02105344719483071, 0.021027834150496361),
    @(-0.015027834314476962, 0.015004680154388161),
    @(-0.0090046803191325964, 0.0089999998279477467),
    @(-0.036112042962649582, 0.036097232238802235),
    @(-0.027097232391928294, 0.027014119460860186),
    @(-0.018014119615438418, 0.018004334748187745),
    @(-0.0090043349029933495, 0.0089999998450052132),
    @(-0.13377498994685055, 0.13309755207294138),
    @(-0.084637976623031186, 0.084127636165046482),
    @(-0.042127636404988955, 0.041999999758678186),
    @(-0.047814913579795615, 0.047749557759448891),
    @(-0.041749557916741509, 0.041671722343597395),
    @(-0.035671722501292358, 0.035427137465199099),
    @(-0.029427137624447042, 0.029272442603957671),
    @(-0.062530072692716487, 0.062168419032403222),
    @(-0.042168419226953358, 0.042020360288905767),
    @(-0.027020360473331451, 0.027001011049820889),
    @(-0.0060010112487773881, 0.0059999998364794749)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
